# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G ("K") values are recalculated; only the G2:G42 numbers change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,7,9,3,8,8,1,6,5,9,5,6,9,8,6,7,6,5,12,4,5,9,7,10,5,7,15,4,6,5,6,3,4,6,5,4,1,7,1,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
